# "Ruoli e Regole - Menu"
# Adds two new menu rows ("Ruoli" and "Utenti") to the Menu sheet, right
# after the existing "Anagrafiche" entry (row 3), pushing the rest of the
# sheet's formatted-but-empty rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a single new row at position 4 - this shifts every row currently
# at 4 or below (the old row 5, 15, 21, 31, 33-36, 39) down by one.
$ws.Rows("4:4").Insert()

# Fill in the new "Ruoli" row (row 4).
# Strings are entered in this order so the shared-string table ends up
# with the same ordering as the target workbook (Utenti, Ruoli, then the
# two new urls).
$ws.Range("B5").Value = "Utenti"
$ws.Range("C5").Value = "Utenti"
$ws.Range("B4").Value = "Ruoli"
$ws.Range("C4").Value = "Ruoli"
$ws.Range("H4").Value = "/registry/list-roles"
$ws.Range("H5").Value = "/registry/list-users"

$ws.Range("A4").Value = 21
$ws.Range("D4").Value = 20
$ws.Range("F4").Value = 1
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0

# Fill in the new "Utenti" row (row 5).
$ws.Range("A5").Value = 22
$ws.Range("D5").Value = 20
$ws.Range("F5").Value = 2
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0

# Leave the selection where the author left it when saving.
$ws.Range("H6").Select()
